$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new bold note in C23 (row 22 intentionally left blank, matching the
# spacer row already present before row 21 -> new content row 23)
$ws.Range("C23").Value = "** Bold outputs show discrepencies between what was predicted and what happened **"

# Match formatting: non-bold black "Aptos Narrow" cell (C3's style) promoted to Bold
$ws.Range("C3").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Font.Bold = $true

# Update selection / view state to match the after-state
$ws.Range("B28").Select() | Out-Null
